$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.898.14"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.810.90"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4283"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07218"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8605"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").Value = "2.034.45"
$ws.Range("E11").Value = "  +12.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.621"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.377"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06895"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008849"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "26.950.16"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "2.281.81"
$ws.Range("E24").Value = "  +12.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.886"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.209"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.889"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08924"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7422"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.411"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05209"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5069"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.747"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1640"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.427"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.241"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.647"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4553"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.801"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
